# Auto-generated script to apply odds updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 15
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("AA2").Value = 29
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 23
$ws.Range("AX2").Value = 4
$ws.Range("AY2").Value = 9

# Row 4
$ws.Range("Q4").Value = 1.85
$ws.Range("R4").Value = 2.05

# Row 5
$ws.Range("G5").Value = 3.3
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3.75
$ws.Range("L5").Value = 2.6
$ws.Range("N5").Value = 12.3
$ws.Range("AA5").Value = 23
$ws.Range("AJ5").Value = 9
$ws.Range("AK5").Value = 19
$ws.Range("AO5").Value = 17
$ws.Range("AS5").Value = 126
$ws.Range("AY5").Value = 11

# Row 7
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 1.36
$ws.Range("O7").Value = 1.11
$ws.Range("P7").Value = 6.5
$ws.Range("Q7").Value = 1.4
$ws.Range("R7").Value = 2.88
$ws.Range("W7").Value = 29
$ws.Range("AJ7").Value = 9
$ws.Range("AN7").Value = 9
$ws.Range("AZ7").Value = 13
$ws.Range("BA7").Value = 15

# Row 8
$ws.Range("G8").Value = 1.22
$ws.Range("H8").Value = 6.25
$ws.Range("J8").Value = 1.62
$ws.Range("K8").Value = 2.88
$ws.Range("L8").Value = 8.5
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 21
$ws.Range("O8").Value = 1.11
$ws.Range("P8").Value = 6.5
$ws.Range("Q8").Value = 1.4
$ws.Range("R8").Value = 2.88
$ws.Range("S8").Value = 1.22
$ws.Range("T8").Value = 4
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("W8").Value = 11
$ws.Range("X8").Value = 7.5
$ws.Range("Z8").Value = 8
$ws.Range("AC8").Value = 21
$ws.Range("AD8").Value = 12
$ws.Range("AH8").Value = 29
$ws.Range("AI8").Value = 51
$ws.Range("AJ8").Value = 29
$ws.Range("AK8").Value = 126
$ws.Range("AN8").Value = 3.5
$ws.Range("AP8").Value = 15
$ws.Range("AQ8").Value = 12
$ws.Range("AT8").Value = 4
$ws.Range("AU8").Value = 9
$ws.Range("AX8").Value = 11
$ws.Range("BA8").Value = 151
$ws.Range("BC8").Value = 251

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 4.1
$ws.Range("J11").Value = 2.2
$ws.Range("K11").Value = 2.25
$ws.Range("L11").Value = 4.35
$ws.Range("M11").Value = 1.02
$ws.Range("N11").Value = 12.5
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 3.85
$ws.Range("Q11").Value = 1.65
$ws.Range("R11").Value = 1.98
$ws.Range("S11").Value = 1.32
$ws.Range("T11").Value = 3.2
$ws.Range("X11").Value = 7.2
$ws.Range("Z11").Value = 10.75
$ws.Range("AA11").Value = 10.5
$ws.Range("AC11").Value = 12
$ws.Range("AD11").Value = 6.6
$ws.Range("AH11").Value = 11
$ws.Range("AI11").Value = 19.5
$ws.Range("AJ11").Value = 11.5
$ws.Range("AK11").Value = 50
$ws.Range("AL11").Value = 29
$ws.Range("AM11").Value = 30
$ws.Range("AN11").Value = 3.6
$ws.Range("AO11").Value = 8
$ws.Range("AP11").Value = 16.5
$ws.Range("AQ11").Value = 26
$ws.Range("AT11").Value = 2.95
$ws.Range("AU11").Value = 7.3
$ws.Range("AX11").Value = 6
$ws.Range("AY11").Value = 22
$ws.Range("AZ11").Value = 27
$ws.Range("BA11").Value = 120

# Row 12
$ws.Range("G12").Value = 1.57
$ws.Range("I12").Value = 4.9
$ws.Range("J12").Value = 2.07
$ws.Range("K12").Value = 2.22
$ws.Range("L12").Value = 5.1
$ws.Range("P12").Value = 3.85
$ws.Range("S12").Value = 1.32
$ws.Range("T12").Value = 3.2
$ws.Range("U12").Value = 1.72
$ws.Range("V12").Value = 2.07
$ws.Range("W12").Value = 6.7
$ws.Range("Y12").Value = 6.8
$ws.Range("Z12").Value = 10.5
$ws.Range("AA12").Value = 10
$ws.Range("AB12").Value = 17.5
$ws.Range("AC12").Value = 11.5
$ws.Range("AE12").Value = 12
$ws.Range("AI12").Value = 24
$ws.Range("AJ12").Value = 13
$ws.Range("AK12").Value = 70
$ws.Range("AL12").Value = 37
$ws.Range("AM12").Value = 35
$ws.Range("AN12").Value = 3.5
$ws.Range("AO12").Value = 7.5
$ws.Range("AP12").Value = 15
$ws.Range("AQ12").Value = 23
$ws.Range("AR12").Value = 45
$ws.Range("AS12").Value = 175
$ws.Range("AT12").Value = 2.95
$ws.Range("AX12").Value = 6.8
$ws.Range("AY12").Value = 29
$ws.Range("AZ12").Value = 32
$ws.Range("BA12").Value = 175
$ws.Range("BB12").Value = 200
$ws.Range("BC12").Value = 400

# Row 14
$ws.Range("G14").Value = 27
$ws.Range("H14").Value = 10.75
$ws.Range("I14").Value = 1.03
$ws.Range("J14").Value = 19.5
$ws.Range("K14").Value = 4.4
$ws.Range("L14").Value = 1.2
$ws.Range("Q14").Value = 1.11
$ws.Range("R14").Value = 5.6
$ws.Range("S14").Value = 1.1
$ws.Range("T14").Value = 5.8
$ws.Range("U14").Value = 2.18
$ws.Range("V14").Value = 1.62
$ws.Range("W14").Value = 200
$ws.Range("X14").Value = 800
$ws.Range("Y14").Value = 200
$ws.Range("AA14").Value = 900
$ws.Range("AB14").Value = 450
$ws.Range("AC14").Value = 50
$ws.Range("AD14").Value = 45
$ws.Range("AE14").Value = 75
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 25
$ws.Range("AI14").Value = 12
$ws.Range("AJ14").Value = 21
$ws.Range("AK14").Value = 8.5
$ws.Range("AL14").Value = 14.5
$ws.Range("AN14").Value = 30
$ws.Range("AO14").Value = 250
$ws.Range("AP14").Value = 120
$ws.Range("AT14").Value = 5.8
$ws.Range("AU14").Value = 14.5
$ws.Range("AX14").Value = 3.85
$ws.Range("AY14").Value = 3.95
$ws.Range("BA14").Value = 5.7
$ws.Range("BB14").Value = 21

# Row 15
$ws.Range("I15").Value = 1.9
$ws.Range("K15").Value = 2.1
$ws.Range("L15").Value = 2.6
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.3
$ws.Range("P15").Value = 3.4
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.75
$ws.Range("S15").Value = 1.44
$ws.Range("T15").Value = 2.63
$ws.Range("AC15").Value = 9.5
$ws.Range("AG15").Value = 301
$ws.Range("AJ15").Value = 9
$ws.Range("AL15").Value = 17
$ws.Range("AM15").Value = 29
$ws.Range("AS15").Value = 251
$ws.Range("AT15").Value = 2.63
$ws.Range("AU15").Value = 8.5

# Row 16
$ws.Range("K16").Value = 1.91
$ws.Range("M16").Value = 1.11
$ws.Range("N16").Value = 6.5
$ws.Range("O16").Value = 1.53
$ws.Range("P16").Value = 2.38
$ws.Range("Q16").Value = 2.7
$ws.Range("R16").Value = 1.44
$ws.Range("U16").Value = 2.2
$ws.Range("V16").Value = 1.62
$ws.Range("W16").Value = 6
$ws.Range("Y16").Value = 10
$ws.Range("Z16").Value = 21
$ws.Range("AE16").Value = 19
$ws.Range("AP16").Value = 29
$ws.Range("AR16").Value = 81
$ws.Range("AS16").Value = 301
$ws.Range("AZ16").Value = 34
